# Generate Report for Handoff
# Updates the localization-status workbook: the first row's source file
# (77879806-...) becomes 1c6ab553-... and has been handed off (status
# "Ready for handoff"), the second row's source file (87db403b-...)
# becomes ffff03f87db4-... The per-language sheets lose their stale
# "Latest Target File"/"Latest Handback File" values (handback hasn't
# happened yet for the new handoff), Content Duplicate flips, and the
# Latest Handoff File/Datetime move forward.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "77879806-9de0-495c-b6c6-b7169e833960"
$oldGuid2 = "87db403b-203f-44e4-a0e0-0276bc3326ca"
$newGuid1 = "1c6ab553-4125-4e21-8fcc-342fe0b6b8b3"
$newGuid2 = "ffff03f87db4-b90d-4a9a-9d75-c5654bdfc1c8"

$status = "Ready for handoff"
$genDate = "2016-09-06 23:17:27"

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e004a855ce08326f1f5a5644f127a3f5f8441729/e2e/"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = $newGuid1 + ".md"
$ov.Range("A3").Value = $newGuid2 + ".md"

$ov.Range("E2").Value = $status
$ov.Range("F2").Value = $status
$ov.Range("E3").Value = $status
$ov.Range("F3").Value = $status

$ov.Range("G2").Value = $genDate
$ov.Range("G3").Value = $genDate

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), $baseUrl + $newGuid1 + ".md", "", "", "e2e\" + $newGuid1 + ".md")
$ov.Hyperlinks.Add($ov.Range("B3"), $baseUrl + $newGuid2 + ".md", "", "", "e2e\" + $newGuid2 + ".md")

$ov.Columns.Item(5).AutoFit()
$ov.Columns.Item(6).AutoFit()

# ---------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de)
# ---------------------------------------------------------------------
$langs = @("zh-cn", "de-de")
foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang)

    $ws.Range("A2").Value = $newGuid1 + ".md"
    $ws.Range("A3").Value = $newGuid2 + ".md"

    $ws.Range("C2").Value = $status
    $ws.Range("C3").Value = $status

    # New handoff xliff (same file re-used for both rows per the new data set)
    $ws.Range("G2").Value = $newGuid1 + ".b19004e21b8bac368638e4e4dc40d70d3dff6a78." + $lang + ".xlf"
    $ws.Range("G3").Value = $newGuid1 + ".b19004e21b8bac368638e4e4dc40d70d3dff6a78." + $lang + ".xlf"

    if ($lang -eq "zh-cn") {
        $ws.Range("H2").Value = "2016-09-06 23:17:22"
        $ws.Range("H3").Value = "2016-09-06 23:17:22"
    } else {
        $ws.Range("H2").Value = $genDate
        $ws.Range("H3").Value = $genDate
    }

    # No handback yet against the new handoff - target/handback file + hyperlinks go away
    $ws.Range("I2").Value = ""
    $ws.Range("J2").Value = ""
    $ws.Range("I3").Value = ""
    $ws.Range("J3").Value = ""
    $ws.Range("I2").ClearFormats()
    $ws.Range("I3").ClearFormats()

    $ws.Range("K2").Value = "0001-01-01 00:00:00"
    $ws.Range("K3").Value = "0001-01-01 00:00:00"

    # Content Duplicate flips: row2 -> False, row3 -> True
    $ws.Range("F2").Value = "False"
    $ws.Range("F3").Value = "True"

    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $baseUrl + $newGuid1 + ".md", "", "", $newGuid1 + ".md")
    $ws.Hyperlinks.Add($ws.Range("A3"), $baseUrl + $newGuid2 + ".md", "", "", $newGuid2 + ".md")

    $ws.Columns.Item(3).AutoFit()
    $ws.Columns.Item(9).AutoFit()
    $ws.Columns.Item(10).AutoFit()
}
